$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.720.29"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "3.381.93"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.08"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.33"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D8").Value = "3.381.91"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.47"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").Value = "3.959.70"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.42"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "3.376.07"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "60.827.39"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.21"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.95"
$ws.Range("E21").Value = "  -6.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.00"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.558"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -5.87%  "
$ws.Range("D27").Value = "3.521.17"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.60"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.93"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.76"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "3.412.44"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.98"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("E40").Value = "  -4.51%  "

# Row 41/42 swap: EnergySwap and Hedera swap positions
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0773"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.64"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.779"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.85"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.66"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("D48").Value = "2.511.02"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.47"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.82"
$ws.Range("E51").Value = "  -1.51%  "